# Auto-update stock values: 2025-12-10 07:56:57 UTC
# Adds a new date column (2025-12-09) to each per-field sheet (sheet index 2-13),
# copying the header/style formatting from the previous last column.
$wb = $excel.ActiveWorkbook

# --- Sheet 2: 시가 (new col 73, width 12) ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(1, 72).Copy()
$ws.Cells.Item(1, 73).PasteSpecial(-4122)
$ws.Cells.Item(1, 73).Value = 20251209
$ws.Application.CutCopyMode = $false
$ws.Columns.Item(73).ColumnWidth = 11.165
$ws.Cells.Item(2, 73).Value = 623.01
$ws.Cells.Item(3, 73).Value = 55.45

# --- Sheet 3: 고가 (new col 73, width 12) ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(1, 72).Copy()
$ws.Cells.Item(1, 73).PasteSpecial(-4122)
$ws.Cells.Item(1, 73).Value = 20251209
$ws.Application.CutCopyMode = $false
$ws.Columns.Item(73).ColumnWidth = 11.165
$ws.Cells.Item(2, 73).Value = 625.87
$ws.Cells.Item(3, 73).Value = 56.22

# --- Sheet 4: 저가 (new col 73, width 12) ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(1, 72).Copy()
$ws.Cells.Item(1, 73).PasteSpecial(-4122)
$ws.Cells.Item(1, 73).Value = 20251209
$ws.Application.CutCopyMode = $false
$ws.Columns.Item(73).ColumnWidth = 11.165
$ws.Cells.Item(2, 73).Value = 621
$ws.Cells.Item(3, 73).Value = 55.09

# --- Sheet 5: 종가 (new col 73, width 12) ---
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(1, 72).Copy()
$ws.Cells.Item(1, 73).PasteSpecial(-4122)
$ws.Cells.Item(1, 73).Value = 20251209
$ws.Application.CutCopyMode = $false
$ws.Columns.Item(73).ColumnWidth = 11.165
$ws.Cells.Item(2, 73).Value = 625.05
$ws.Cells.Item(3, 73).Value = 56.01

# --- Sheet 6: 거래량 (new col 73, width 12) ---
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(1, 72).Copy()
$ws.Cells.Item(1, 73).PasteSpecial(-4122)
$ws.Cells.Item(1, 73).Value = 20251209
$ws.Application.CutCopyMode = $false
$ws.Columns.Item(73).ColumnWidth = 11.165
$ws.Cells.Item(2, 73).Value = 37155999
$ws.Cells.Item(3, 73).Value = 48735976

# --- Sheet 7: s20 (new col 54, width 10) ---
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(1, 53).Copy()
$ws.Cells.Item(1, 54).PasteSpecial(-4122)
$ws.Cells.Item(1, 54).Value = 20251209
$ws.Application.CutCopyMode = $false
$ws.Columns.Item(54).ColumnWidth = 9.165
$ws.Cells.Item(2, 54).Value = 99
$ws.Cells.Item(3, 54).Value = 15

# --- Sheet 8: s60 (new col 14, width 10) ---
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(1, 13).Copy()
$ws.Cells.Item(1, 14).PasteSpecial(-4122)
$ws.Cells.Item(1, 14).Value = 20251209
$ws.Application.CutCopyMode = $false
$ws.Columns.Item(14).ColumnWidth = 9.165
$ws.Cells.Item(2, 14).Value = 79
$ws.Cells.Item(3, 14).Value = 13

# --- Sheet 9: z20 (new col 54, width 10) ---
$ws = $wb.Worksheets.Item(9)
$ws.Cells.Item(1, 53).Copy()
$ws.Cells.Item(1, 54).PasteSpecial(-4122)
$ws.Cells.Item(1, 54).Value = 20251209
$ws.Application.CutCopyMode = $false
$ws.Columns.Item(54).ColumnWidth = 9.165
$ws.Cells.Item(2, 54).Value = 53
$ws.Cells.Item(3, 54).Value = -30

# --- Sheet 10: z60 (new col 14, width 10) ---
$ws = $wb.Worksheets.Item(10)
$ws.Cells.Item(1, 13).Copy()
$ws.Cells.Item(1, 14).PasteSpecial(-4122)
$ws.Cells.Item(1, 14).Value = 20251209
$ws.Application.CutCopyMode = $false
$ws.Columns.Item(14).ColumnWidth = 9.165
$ws.Cells.Item(2, 14).Value = 62
$ws.Cells.Item(3, 14).Value = -86

# --- Sheet 11: gap (new col 54, width 12) ---
$ws = $wb.Worksheets.Item(11)
$ws.Cells.Item(1, 54).NumberFormat = "@"
$ws.Cells.Item(1, 54).Value = "20251209"
$ws.Cells.Item(1, 53).Copy()
$ws.Cells.Item(1, 54).PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Columns.Item(54).ColumnWidth = 11.165
$ws.Cells.Item(2, 54).Value = 102
$ws.Cells.Item(3, 54).Value = 79

# --- Sheet 12: std (new col 35, width 12) ---
$ws = $wb.Worksheets.Item(12)
$ws.Cells.Item(1, 35).NumberFormat = "@"
$ws.Cells.Item(1, 35).Value = "20251209"
$ws.Cells.Item(1, 34).Copy()
$ws.Cells.Item(1, 35).PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Columns.Item(35).ColumnWidth = 11.165
$ws.Cells.Item(2, 35).Value = 3.6
$ws.Cells.Item(3, 35).Value = 37.5

# --- Sheet 13: quant (new col 14, width 12) ---
$ws = $wb.Worksheets.Item(13)
$ws.Cells.Item(1, 14).NumberFormat = "@"
$ws.Cells.Item(1, 14).Value = "20251209"
$ws.Cells.Item(1, 13).Copy()
$ws.Cells.Item(1, 14).PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Columns.Item(14).ColumnWidth = 11.165
$ws.Cells.Item(2, 14).Value = 32
$ws.Cells.Item(3, 14).Value = 36
